# Applies the "Added calculation of full DNS history (v0.14)" change:
#  - Renames the dns1/dns2/dns3 variable rows to tdns1/tdns2/tdns3
#  - Fills in the fullname/category/source/units/freq/st/d1/d2 metadata
#    for those three rows on the "baseline-variables" sheet
#  - Updates the sheet's active selection to I4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("baseline-variables")

# Rename the variable id column first (in-place edits of the existing
# dns1/dns2/dns3 shared strings -> tdns1/tdns2/tdns3), matching the
# order the new shared-string table entries were produced in.
$ws.Range("A23").Value = "tdns1"
$ws.Range("A24").Value = "tdns2"
$ws.Range("A25").Value = "tdns3"

# Row 23 - tdns1
$ws.Range("B23").Value = "Treasury-FFR Spread Level (10-Year Level)"
$ws.Range("C23").Value = "Interest Rates"
$ws.Range("D23").Value = "calc"
$ws.Range("F23").Value = "%"
$ws.Range("G23").Value = "m"
$ws.Range("I23").Value = "d"
$ws.Range("J23").Value = "base"
$ws.Range("K23").Value = "none"

# Row 24 - tdns2
$ws.Range("B24").Value = "Treasury-FFR Spread Slope (Negative of 10Y-3M Spread)"
$ws.Range("C24").Value = "Interest Rates"
$ws.Range("D24").Value = "calc"
$ws.Range("F24").Value = "%"
$ws.Range("G24").Value = "m"
$ws.Range("I24").Value = "d"
$ws.Range("J24").Value = "base"
$ws.Range("K24").Value = "none"

# Row 25 - tdns3
$ws.Range("B25").Value = "Treasury-FFR Spread Curvature"
$ws.Range("C25").Value = "Interest Rates"
$ws.Range("D25").Value = "calc"
$ws.Range("F25").Value = "%"
$ws.Range("G25").Value = "m"
$ws.Range("I25").Value = "d"
$ws.Range("J25").Value = "base"
$ws.Range("K25").Value = "none"

# Update active selection on the sheet
$ws.Activate()
$ws.Range("I4").Select()
